# "Аналіз проектних ризиків.xlsx" - sheet "Ризики" lists project risks with
# a running "#" number in column A (rows 2-16). The rows were renumbered so
# column A now simply counts 1..15 top to bottom: row 2 already holds 1 and
# stays as-is, while rows 3-16 are rewritten to the sequential values 2-15.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNumbers = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15)

for ($i = 0; $i -lt $newNumbers.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 1).Value = $newNumbers[$i]
}
